# Croatia HNL workbook update - 21-04-2024 14:32
# Two new match rows are inserted (2024-04-20 games) before the existing
# row that used to be row 152 (NK Lokomotiva Zagreb vs Dinamo Zagreb,
# 2024-04-21). This shifts the four existing trailing rows (old rows
# 152-155) down to rows 154-157, and the two new rows become rows 152-153.
#
# We avoid Rows.Insert() because that causes Excel to synthesize brand
# new cell styles (it cannot always reuse the existing border-less bold
# style for the newly inserted blank row), which would needlessly grow
# styles.xml. Instead we shift the data "manually": copy formats with
# PasteSpecial (which reuses the existing style indexes) and then copy
# values on top, working from the bottom row upwards so we never
# overwrite data before it has been read.
#
# Note: column A ("id") is a plain row-sequence counter (A<n> always
# equals n-2) and is NOT part of the shifted data - it keeps whatever
# value already belongs to its row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Shift-Row([int]$fromRow, [int]$toRow) {
    $srcFull = $ws.Range("B" + $fromRow + ":AC" + $fromRow)
    $dstFull = $ws.Range("B" + $toRow + ":AC" + $toRow)

    # Copy cell formatting (styles) first so the destination row ends up
    # with the exact same style indexes as the source row (date column
    # custom date format, etc.)
    $srcFull.Copy()
    $dstFull.PasteSpecial(-4122)

    # Now copy the actual values across (array copy keeps blank cells
    # blank instead of writing zeros).
    $vals = $srcFull.Value()
    $dstFull.Value = $vals
}

# Shift existing rows 152-155 down to 154-157, starting from the bottom
# so source data is not clobbered before being read.
Shift-Row 155 157
Shift-Row 154 156
Shift-Row 153 155
Shift-Row 152 154

# Rows 156-157 are brand new rows (the sheet previously ended at row
# 155), so - unlike rows 154-155, which already owned their id-column
# cell - they need column A populated explicitly (still following the
# row-sequence rule A<n> = n-2, same style "1" as every other id cell).
$ws.Range("A151").Copy()
$ws.Range("A156:A157").PasteSpecial(-4122)
$ws.Range("A156").Value = 154
$ws.Range("A157").Value = 155

# Rows 152-153 already carry the correct id-column (A) value/style
# (150/151, both style "1") untouched from before the edit, since that
# column is a plain row-sequence counter. Only the date column (E) needs
# its format (custom date numFmt) applied before writing the new dates.
$ws.Range("E151").Copy()
$ws.Range("E152:E153").PasteSpecial(-4122)

# New row 152: NK Osijek vs Istra 1961, 2024-04-20 12:00
$ws.Range("A152").Value = 150
$ws.Range("B152").Value = 6962505
$ws.Range("C152").Value = "Croatia HNL"
$ws.Range("D152").Value = "Croatia HNL"
$ws.Range("E152").Value = 45402.5
$ws.Range("F152").Value = "NK Osijek"
$ws.Range("G152").Value = "Istra 1961"
$ws.Range("H152").Value = 1
$ws.Range("I152").Value = 2
$ws.Range("J152").Value = "A"
$ws.Range("K152").Value = 1.55
$ws.Range("L152").Value = 3.8
$ws.Range("M152").Value = 6.5
$ws.Range("N152").Value = 1.615
$ws.Range("O152").Value = 3.6
$ws.Range("P152").Value = 6
$ws.Range("Q152").Value = -0.75
$ws.Range("R152").Value = 1.8
$ws.Range("S152").Value = 2.05
$ws.Range("T152").Value = 2.25
$ws.Range("U152").Value = 1.875
$ws.Range("V152").Value = 1.975
$ws.Range("W152").Value = -1
$ws.Range("X152").Value = -1
$ws.Range("Y152").Value = 5
$ws.Range("Z152").Value = -1
$ws.Range("AA152").Value = 1.05
$ws.Range("AB152").Value = 0.875
$ws.Range("AC152").Value = -1

# New row 153: Slaven Belupo vs Hajduk Split, 2024-04-20 14:10
$ws.Range("A153").Value = 151
$ws.Range("B153").Value = 6954858
$ws.Range("C153").Value = "Croatia HNL"
$ws.Range("D153").Value = "Croatia HNL"
$ws.Range("E153").Value = 45402.59027777778
$ws.Range("F153").Value = "Slaven Belupo"
$ws.Range("G153").Value = "Hajduk Split"
$ws.Range("H153").Value = 0
$ws.Range("I153").Value = 1
$ws.Range("J153").Value = "A"
$ws.Range("K153").Value = 5.75
$ws.Range("L153").Value = 4
$ws.Range("M153").Value = 1.571
$ws.Range("N153").Value = 5
$ws.Range("O153").Value = 3.5
$ws.Range("P153").Value = 1.727
$ws.Range("Q153").Value = 0.75
$ws.Range("R153").Value = 1.875
$ws.Range("S153").Value = 1.975
$ws.Range("T153").Value = 2.5
$ws.Range("U153").Value = 1.925
$ws.Range("V153").Value = 1.925
$ws.Range("W153").Value = -1
$ws.Range("X153").Value = -1
$ws.Range("Y153").Value = 0.7270000000000001
$ws.Range("Z153").Value = -0.5
$ws.Range("AA153").Value = 0.4875
$ws.Range("AB153").Value = -1
$ws.Range("AC153").Value = 0.925
